$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-09-18T19:07:03"
$ws.Range("W4").Value = 122.75
$ws.Range("X4").Value = 115
$ws.Range("Z4").Value = 39.35
$ws.Range("W6").Value = -3.73
$ws.Range("X6").Value = -2.34
$ws.Range("Y6").Value = -0.98
$ws.Range("Z6").Value = -1.38
$ws.Range("Y7").Value = 85.87
$ws.Range("W8").Value = 29.42
$ws.Range("X8").Value = 39.42
$ws.Range("W9").Value = 127.46
$ws.Range("X9").Value = 119.59
$ws.Range("Y9").Value = 41.77
$ws.Range("Z9").Value = 41.86
$ws.Range("W11").Value = 0.98
$ws.Range("X11").Value = 2.24
$ws.Range("Y11").Value = 1.46
$ws.Range("Z11").Value = 1.13
$ws.Range("W13").Value = 29.42
$ws.Range("X13").Value = 39.42
$ws.Range("W14").Value = 127.56
$ws.Range("X14").Value = 140.56
$ws.Range("Y14").Value = 173.42
$ws.Range("Z14").Value = 41.9
$ws.Range("W15").Value = 0
$ws.Range("X15").Value = 20.97
$ws.Range("Y15").Value = 131.6
$ws.Range("W16").Value = 1.08
$ws.Range("X16").Value = 2.24
$ws.Range("Y16").Value = 1.51
$ws.Range("Z16").Value = 1.17
$ws.Range("W18").Value = 29.42
$ws.Range("X18").Value = 39.42
$ws.Range("W19").Value = 122.66
$ws.Range("X19").Value = 115
$ws.Range("Y19").Value = 39.29
$ws.Range("Z19").Value = 39.31
$ws.Range("W21").Value = -3.82
$ws.Range("X21").Value = -2.34
$ws.Range("Y21").Value = -1.02
$ws.Range("Z21").Value = -1.42
$ws.Range("W23").Value = 29.42
$ws.Range("X23").Value = 39.42
$ws.Range("W24").Value = 122.66
$ws.Range("X24").Value = 115
$ws.Range("Z24").Value = 39.31
$ws.Range("W26").Value = -3.82
$ws.Range("X26").Value = -2.34
$ws.Range("Y26").Value = -1.02
$ws.Range("Z26").Value = -1.42
$ws.Range("Y27").Value = 85.87
$ws.Range("W28").Value = 29.42
$ws.Range("X28").Value = 39.42
$ws.Range("W29").Value = 122.21
$ws.Range("X29").Value = 114.71
$ws.Range("Y29").Value = 39.02
$ws.Range("Z29").Value = 39.16
$ws.Range("W31").Value = -4.27
$ws.Range("X31").Value = -2.64
$ws.Range("Y31").Value = -1.29
$ws.Range("Z31").Value = -1.57
$ws.Range("W33").Value = 29.42
$ws.Range("X33").Value = 39.42
$ws.Range("X34").Value = 144
$ws.Range("Y34").Value = 175.32
$ws.Range("Z34").Value = 43.51
$ws.Range("W35").Value = 0
$ws.Range("X35").Value = 20.97
$ws.Range("Y35").Value = 131.6
$ws.Range("W36").Value = 4.89
$ws.Range("X36").Value = 5.69
$ws.Range("Y36").Value = 3.41
$ws.Range("Z36").Value = 2.78
$ws.Range("W37").Value = 13.63
$ws.Range("W38").Value = 29.42
$ws.Range("X38").Value = 39.42
$ws.Range("W39").Value = 122.75
$ws.Range("X39").Value = 115
$ws.Range("Z39").Value = 39.35
$ws.Range("W41").Value = -3.73
$ws.Range("X41").Value = -2.34
$ws.Range("Y41").Value = -0.98
$ws.Range("Z41").Value = -1.38
$ws.Range("Y42").Value = 85.87
$ws.Range("W43").Value = 29.42
$ws.Range("X43").Value = 39.42
$ws.Range("W44").Value = 125.05
$ws.Range("X44").Value = 116.65
$ws.Range("Y44").Value = 40.19
$ws.Range("Z44").Value = 40.32
$ws.Range("W46").Value = -1.43
$ws.Range("X46").Value = -0.7
$ws.Range("Y46").Value = -0.12
$ws.Range("Z46").Value = -0.4
$ws.Range("W48").Value = 29.42
$ws.Range("X48").Value = 39.42
$ws.Range("W49").Value = 131.59
$ws.Range("X49").Value = 121.44
$ws.Range("Y49").Value = 42.12
$ws.Range("Z49").Value = 42.42
$ws.Range("W51").Value = 5.11
$ws.Range("X51").Value = 4.1
$ws.Range("Y51").Value = 1.81
$ws.Range("Z51").Value = 1.7
$ws.Range("W53").Value = 29.42
$ws.Range("X53").Value = 39.42
$ws.Range("W54").Value = 126.09
$ws.Range("X54").Value = 118.05
$ws.Range("Y54").Value = 41.43
$ws.Range("Z54").Value = 41.73
$ws.Range("W56").Value = -0.39
$ws.Range("X56").Value = 0.71
$ws.Range("Z56").Value = 1
$ws.Range("W58").Value = 29.43
$ws.Range("X58").Value = 39.42
$ws.Range("W59").Value = 128.26
$ws.Range("X59").Value = 119.26
$ws.Range("Y59").Value = 41.39
$ws.Range("Z59").Value = 41.6
$ws.Range("W61").Value = 1.78
$ws.Range("X61").Value = 1.92
$ws.Range("Y61").Value = 1.08
$ws.Range("Z61").Value = 0.87
$ws.Range("W63").Value = 29.42
$ws.Range("X63").Value = 39.42
$ws.Range("W64").Value = 130
$ws.Range("X64").Value = 120.51
$ws.Range("Y64").Value = 41.95
$ws.Range("Z64").Value = 42.2
$ws.Range("W66").Value = 3.52
$ws.Range("X66").Value = 3.16
$ws.Range("Y66").Value = 1.64
$ws.Range("Z66").Value = 1.48
$ws.Range("W68").Value = 29.42
$ws.Range("X68").Value = 39.42
$ws.Range("X69").Value = 121.19
$ws.Range("Z69").Value = 42.6
$ws.Range("W71").Value = 3.52
$ws.Range("X71").Value = 3.84
$ws.Range("Y71").Value = 2.03
$ws.Range("Z71").Value = 1.87
$ws.Range("W73").Value = 29.42
$ws.Range("X73").Value = 39.42
$ws.Range("W74").Value = 128.56
$ws.Range("X74").Value = 119.42
$ws.Range("Y74").Value = 41.39
$ws.Range("Z74").Value = 41.64
$ws.Range("W76").Value = 2.08
$ws.Range("X76").Value = 2.08
$ws.Range("Y76").Value = 1.08
$ws.Range("Z76").Value = 0.92
$ws.Range("W78").Value = 29.43
$ws.Range("X78").Value = 39.42
$ws.Range("W79").Value = 129.16
$ws.Range("X79").Value = 119.78
$ws.Range("Y79").Value = 41.62
$ws.Range("Z79").Value = 41.87
$ws.Range("W81").Value = 2.68
$ws.Range("X81").Value = 2.44
$ws.Range("Z81").Value = 1.15
$ws.Range("W83").Value = 29.42
$ws.Range("X83").Value = 39.42
$ws.Range("W84").Value = 126.29
$ws.Range("X84").Value = 117.74
$ws.Range("Y84").Value = 41.77
$ws.Range("Z84").Value = 42.25
$ws.Range("W86").Value = -0.19
$ws.Range("X86").Value = 0.39
$ws.Range("Z86").Value = 1.52
$ws.Range("W88").Value = 29.42
$ws.Range("X88").Value = 39.42
$ws.Range("W89").Value = 122.21
$ws.Range("X89").Value = 114.71
$ws.Range("Y89").Value = 39.02
$ws.Range("Z89").Value = 39.16
$ws.Range("W91").Value = -4.27
$ws.Range("X91").Value = -2.64
$ws.Range("Y91").Value = -1.29
$ws.Range("Z91").Value = -1.57
$ws.Range("W93").Value = 29.42
$ws.Range("X93").Value = 39.42
